$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: insert new "Condition" column with custom width ---
$ws.Columns.Item(3).ColumnWidth = 24.88671875

# Grab the two distinct D-column number formats that already exist in the
# workbook before we start overwriting any of the source cells:
#   - D3 (style s=3): plain percent format, default font (no explicit color)
#   - D4 (style s=4): plain percent format, font with no explicit color attr
# Re-applying (copying) them onto the new D column cells makes the engine
# map every cell back onto the pre-existing style indexes instead of
# minting brand new, duplicate ones. Grab D4's original (distinct) format
# first, before it gets overwritten by the D3 bulk-format paste below.
$ws.Range("D4").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D3:D8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 1 (header row) ---
$ws.Range("A1").Value = "Machine Learning Models"
$ws.Range("B1").Value = "Types"
$ws.Range("C1").Value = "Condition"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").WrapText = $true
$ws.Range("D1").Value = "Accurracy (target > 75%)"
$ws.Range("E1").Value = "Comments"

# --- Row 2 ---
$ws.Range("A2").Value = "Logistic Regression "
$ws.Range("B2").Value = "Supervised"
$ws.Range("C2").Value = "Original Imbalanced Data"
$ws.Range("C2").WrapText = $true
$ws.Range("D2").Value = 0.95
$ws.Range("E2").Value = "No winning team data, class 0 has 823, class 1 has 31"

# --- Row 3 ---
$ws.Range("A3").Value = "Logistic Regression "
$ws.Range("B3").Value = "Supervised"
$ws.Range("C3").Value = "Random Over Sampler to Balance Data"
$ws.Range("C3").WrapText = $true
$ws.Range("D3").Value = 0.88
$ws.Range("E3").Value = "Improvement on winning team.  Class 0 and Class 1 have 614"
$ws.Rows.Item(3).RowHeight = 28.8

# --- Row 4 ---
$ws.Range("A4").Value = "Decision Tree Classifier"
$ws.Range("B4").Value = "Supervised"
$ws.Range("C4").Value = "Random Over Sampler to Balance Data"
$ws.Range("C4").WrapText = $true
$ws.Range("D4").Value = 0.98
$ws.Range("E4").Value = "Consistency in results for both winning teams and losing teams"
$ws.Rows.Item(4).RowHeight = 28.8

# --- Row 5 ---
$ws.Range("A5").Value = "Decision Tree Classifier "
$ws.Range("B5").Value = "Supervised"
$ws.Range("C5").Value = "Random Over Sampler to Balance Data With Class Weight"
$ws.Range("C5").WrapText = $true
$ws.Range("D5").Value = 0.99
$ws.Range("E5").Value = "Consistency in results for both winning teams and losing teams"
$ws.Rows.Item(5).RowHeight = 43.2

# --- Row 6 ---
$ws.Range("A6").Value = "Decision Tree Classifier"
$ws.Range("B6").Value = "Supervised"
$ws.Range("C6").Value = "Under Sampling With Class Weight"
$ws.Range("C6").WrapText = $true
$ws.Range("D6").Value = 0.77
$ws.Range("E6").Value = "Equal samples for Class 0 and Class 1 to 31"
$ws.Rows.Item(6).RowHeight = 28.8

# --- Row 7 (new row) ---
$ws.Range("A7").Value = "Decision Tree Classifier"
$ws.Range("B7").Value = "Supervised"
$ws.Range("C7").Value = "Over Sampling With Class Weight"
$ws.Range("C7").WrapText = $true
$ws.Range("D7").Value = 0.97
$ws.Range("E7").Value = "Equal samples for Class 0 and Class 1 to 823"
$ws.Rows.Item(7).RowHeight = 28.8

# --- Row 8 (new row) ---
$ws.Range("A8").Value = "Decision Tree Classifier"
$ws.Range("B8").Value = "Supervised"
$ws.Range("C8").Value = "SMOTE (Synthetic Minority Over-sampling Technique) With Class Weight"
$ws.Range("C8").WrapText = $true
$ws.Range("D8").Value = 0.97
$ws.Range("E8").Value = "Equal samples for Class 0 and Class 1 to 823"
$ws.Rows.Item(8).RowHeight = 43.2

# --- Column E: clear the old custom width so it reverts to default ---
$ws.Columns.Item(5).ColumnWidth = 8.43

# --- Selection matches the final saved state ---
$ws.Range("E9").Select()

Write-Host "edit complete"
